$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1: "This text is bold But this isn't"  ->  "This is a bold header"
# New text is entirely bold; paragraph loses its line-spacing setting but
# keeps center alignment; the gramStart/gramEnd proofErr markers disappear
# along with the "But" / " this isn't" runs.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1Range = $p1.Range

$p1Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">This </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>is</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>a bold header</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$p1Range.InsertXML($p1Xml) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 2: "This isn't bold too" -> "This paragraph is a body paragraph"
# New paragraph gets spacing/justify/font/size changes. Because this is the
# last paragraph in the body, its Range carries the document's final
# paragraph mark, and InsertXML on that range inserts a sibling rather than
# replacing the mark's own stored pPr/rPr. Work around it by appending a
# throwaway trailing paragraph first (so paragraph 2 is no longer last),
# doing the full replace, then deleting the throwaway paragraph again.
# ---------------------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null

$p2 = $d.Paragraphs(2)
$p2Range = $p2.Range

$p2Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:spacing w:line="360" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>This paragraph is a body paragraph</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$p2Range.InsertXML($p2Xml) | Out-Null

# Remove the throwaway trailing paragraph we added above.
$trailing = $d.Paragraphs($d.Paragraphs.Count)
$trailing.Range.Delete() | Out-Null
